# Apply updated "dSF" (column F) values as per repull/recalculation of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F3").Value = 5
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = -1
$ws.Range("F8").Value = 5
$ws.Range("F9").Value = -4
$ws.Range("F10").Value = -4
$ws.Range("F11").Value = -1
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = -1
$ws.Range("F17").Value = 8
$ws.Range("F18").Value = -6
$ws.Range("F19").Value = 1
$ws.Range("F20").Value = -1
$ws.Range("F21").Value = 1
